# Revert "weekly training updates"
# This reverts the earlier commit that (a) pushed the weekly schedule's
# start time forward to 10:00 AM, (b) shifted the SUN..SAT day headers one
# column to the right (wrapping SUN to the end), (c) filled in a bunch of
# per-day/per-slot training-log notes, and (d) truncated the auto-filled
# time column partway down the sheet. We restore the prior state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Schedule")

# 1. Schedule start time (E3) back to 5:00 AM (was pushed to 10:00 AM).
$ws.Range("E3").Value = 0.20833333333333334

# 2. Day-of-week header row back to SUN, MON, TUES, WED, THURS, FRI, SAT
#    in columns C:I (previously rotated so SUN had been moved to the end).
$ws.Range("C5").Value = "SUN"
$ws.Range("D5").Value = "MON"
$ws.Range("E5").Value = "TUES"
$ws.Range("F5").Value = "WED"
$ws.Range("G5").Value = "THURS"
$ws.Range("H5").Value = "FRI"
$ws.Range("I5").Value = "SAT"

# 3. Clear the training-log notes that had been typed into the grid.
$notesCells = @(
    "C7","D7","E7","F7","G7",
    "C8","D8","E8","F8","G8",
    "C10","D10","F10","G10",
    "C15","D15","E15","F15","G15",
    "C16","D16","E16","F16","G16",
    "C17","D17","E17","F17","G17",
    "C18","D18",
    "D19","G19",
    "D20","G20"
)
foreach ($addr in $notesCells) {
    $ws.Range($addr).ClearContents()
}

# 4. Restore the auto-incrementing time formula for rows 21:32 (the shared
#    formula previously only spanned B8:B20; it should span B8:B32).
for ($r = 21; $r -le 32; $r++) {
    $prevAddr = "B" + ($r - 1)
    $ws.Cells.Item($r, 2).Formula = "=" + $prevAddr + "+TIME(0,Interval,0)"
}

$wb.Save()
